$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Helper: read a full data row (columns B..AC) from a given row number
# ---------------------------------------------------------------------------
function Get-RowValues($ws, $row) {
    $vals = @()
    for ($c = 2; $c -le 29; $c++) {
        $vals += , $ws.Cells.Item($row, $c).Value()
    }
    return $vals
}

function Set-RowValues($ws, $row, $vals) {
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------------
# 1) Swap data (columns B..AC) between rows 95 and 96
# ---------------------------------------------------------------------------
$row95 = Get-RowValues $ws 95
$row96 = Get-RowValues $ws 96
Set-RowValues $ws 95 $row96
Set-RowValues $ws 96 $row95

# ---------------------------------------------------------------------------
# 2) Rotate data (columns B..AC) among rows 104, 105, 106, 107
#    new104 = old107 ; new105 = old106 ; new106 = old105 ; new107 = old104
# ---------------------------------------------------------------------------
$row104 = Get-RowValues $ws 104
$row105 = Get-RowValues $ws 105
$row106 = Get-RowValues $ws 106
$row107 = Get-RowValues $ws 107

Set-RowValues $ws 104 $row107
Set-RowValues $ws 105 $row106
Set-RowValues $ws 106 $row105
Set-RowValues $ws 107 $row104

# ---------------------------------------------------------------------------
# 3) Append three new rows (119, 120, 121)
# ---------------------------------------------------------------------------

# Row 119
$ws.Cells.Item(119, 1).Value = 117
$ws.Cells.Item(119, 2).Value = 7719643
$ws.Cells.Item(119, 3).Value = "Estonia Meistriliiga"
$ws.Cells.Item(119, 4).Value = "Estonia Meistriliiga"
$ws.Cells.Item(119, 5).Value = 45367.3125
$ws.Cells.Item(119, 6).Value = "JK Tallinna Kalev"
$ws.Cells.Item(119, 7).Value = "FC Levadia Tallinn"
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 9).Value = 2
$ws.Cells.Item(119, 10).Value = "A"
$ws.Cells.Item(119, 11).Value = 5
$ws.Cells.Item(119, 12).Value = 4.5
$ws.Cells.Item(119, 13).Value = 1.444
$ws.Cells.Item(119, 14).Value = 8
$ws.Cells.Item(119, 15).Value = 6.5
$ws.Cells.Item(119, 16).Value = 1.2
$ws.Cells.Item(119, 17).Value = 1.75
$ws.Cells.Item(119, 18).Value = 1.975
$ws.Cells.Item(119, 19).Value = 1.825
$ws.Cells.Item(119, 20).Value = 3
$ws.Cells.Item(119, 21).Value = 1.825
$ws.Cells.Item(119, 22).Value = 1.975
$ws.Cells.Item(119, 23).Value = -1
$ws.Cells.Item(119, 24).Value = -1
$ws.Cells.Item(119, 25).Value = 0.2
$ws.Cells.Item(119, 26).Value = -0.5
$ws.Cells.Item(119, 27).Value = 0.4125
$ws.Cells.Item(119, 28).Value = -1
$ws.Cells.Item(119, 29).Value = 0.9750000000000001

# Row 120
$ws.Cells.Item(120, 1).Value = 118
$ws.Cells.Item(120, 2).Value = 7721007
$ws.Cells.Item(120, 3).Value = "Estonia Meistriliiga"
$ws.Cells.Item(120, 4).Value = "Estonia Meistriliiga"
$ws.Cells.Item(120, 5).Value = 45367.39583333334
$ws.Cells.Item(120, 6).Value = "JK Trans Narva"
$ws.Cells.Item(120, 7).Value = "JK Tammeka Tartu"
$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(120, 9).Value = 5
$ws.Cells.Item(120, 10).Value = "A"
$ws.Cells.Item(120, 11).Value = 2.25
$ws.Cells.Item(120, 12).Value = 3.3
$ws.Cells.Item(120, 13).Value = 2.75
$ws.Cells.Item(120, 14).Value = 2.1
$ws.Cells.Item(120, 15).Value = 3.25
$ws.Cells.Item(120, 16).Value = 3
$ws.Cells.Item(120, 17).Value = -0.25
$ws.Cells.Item(120, 18).Value = 1.875
$ws.Cells.Item(120, 19).Value = 1.925
$ws.Cells.Item(120, 20).Value = 2.5
$ws.Cells.Item(120, 21).Value = 1.825
$ws.Cells.Item(120, 22).Value = 1.975
$ws.Cells.Item(120, 23).Value = -1
$ws.Cells.Item(120, 24).Value = -1
$ws.Cells.Item(120, 25).Value = 2
$ws.Cells.Item(120, 26).Value = -1
$ws.Cells.Item(120, 27).Value = 0.925
$ws.Cells.Item(120, 28).Value = 0.825
$ws.Cells.Item(120, 29).Value = -1

# Row 121
$ws.Cells.Item(121, 1).Value = 119
$ws.Cells.Item(121, 2).Value = 7721087
$ws.Cells.Item(121, 3).Value = "Estonia Meistriliiga"
$ws.Cells.Item(121, 4).Value = "Estonia Meistriliiga"
$ws.Cells.Item(121, 5).Value = 45367.39583333334
$ws.Cells.Item(121, 6).Value = "Paide Linnameeskond"
$ws.Cells.Item(121, 7).Value = "FC Flora Tallinn"
$ws.Cells.Item(121, 8).Value = 2
$ws.Cells.Item(121, 9).Value = 1
$ws.Cells.Item(121, 10).Value = "H"
$ws.Cells.Item(121, 11).Value = 2.2
$ws.Cells.Item(121, 12).Value = 3.3
$ws.Cells.Item(121, 13).Value = 2.8
$ws.Cells.Item(121, 14).Value = 1.85
$ws.Cells.Item(121, 15).Value = 3.6
$ws.Cells.Item(121, 16).Value = 3.4
$ws.Cells.Item(121, 17).Value = -0.5
$ws.Cells.Item(121, 18).Value = 1.9
$ws.Cells.Item(121, 19).Value = 1.9
$ws.Cells.Item(121, 20).Value = 2.5
$ws.Cells.Item(121, 21).Value = 1.95
$ws.Cells.Item(121, 22).Value = 1.85
$ws.Cells.Item(121, 23).Value = 0.8500000000000001
$ws.Cells.Item(121, 24).Value = -1
$ws.Cells.Item(121, 25).Value = -1
$ws.Cells.Item(121, 26).Value = 0.8999999999999999
$ws.Cells.Item(121, 27).Value = -1
$ws.Cells.Item(121, 28).Value = 0.95
$ws.Cells.Item(121, 29).Value = -1

# Apply the same formatting as the other "A" column cells (bold/border,
# centred) and the date-format style used for column E, matching the rest of
# the sheet, by copying formats only (values already written above are kept).
$ws.Range("A118").Copy()
$ws.Range("A119:A121").PasteSpecial(-4122)
$ws.Range("E118").Copy()
$ws.Range("E119:E121").PasteSpecial(-4122)
$excel.CutCopyMode = $false
